$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("initial")
$ws2 = $wb.Worksheets.Item("line_imp")

# --- Sheet "initial": add q_lim column, update bus data, swap pv/pq types ---
$ws1.Range("G1").Value = "q_lim"

# Bus 1 (slack): V goes from 1 to 1.05; clear P and Q (no longer specified for slack)
$ws1.Range("C2").Value = 1.05
$ws1.Range("E2").ClearContents()
$ws1.Range("F2").ClearContents()

# Bus 2: switches from "pv" to "pq"; clear V/T; set new P/Q
$ws1.Range("B3").Value = "pq"
$ws1.Range("C3").ClearContents()
$ws1.Range("D3").ClearContents()
$ws1.Range("E3").Value = -4
$ws1.Range("F3").Value = -2.5

# Bus 3: switches from "pq" to "pv"; set V; clear T; set new P
$ws1.Range("B4").Value = "pv"
$ws1.Range("C4").Value = 1.04
$ws1.Range("D4").ClearContents()
$ws1.Range("E4").Value = 2

# --- Sheet "line_imp": update R/X values for the three lines ---
$ws2.Range("B2").Value = 0.02
$ws2.Range("C2").Value = 0.04

$ws2.Range("B3").Value = 0.01
$ws2.Range("C3").Value = 0.03

$ws2.Range("B4").Value = 0.0125
$ws2.Range("C4").Value = 0.025

# --- View state: active sheet moves from "line_imp" to "initial", selections update ---
[void]$ws2.Activate()
[void]$ws2.Range("F1:G1").Select()
[void]$ws1.Activate()
[void]$ws1.Range("G2").Select()
